# Swap the per-observation values between row 2 and row 3 for the
# columns that actually differ between the two records:
#   A (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
#   G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
# All other columns already hold identical values on both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $value2 = $cell2.Value2
    $value3 = $cell3.Value2

    $cell2.Value2 = $value3
    $cell3.Value2 = $value2
}
